$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Part 1: merge the "/bitnami/" run sequence into a single run and drop
# the spell-check proofErr markers around it.
# -----------------------------------------------------------------------
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Execute(
    "/bitnami/apache2/ et non /",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "/bitnami/apache2/ et non /", 2
) | Out-Null

# -----------------------------------------------------------------------
# Part 2: remove the (hidden) _GoBack bookmark that currently sits in the
# empty paragraph right after "29/05/2018".
# -----------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# -----------------------------------------------------------------------
# Part 3: append the new entries for 12/06/2018 at the end of the report,
# then re-create the _GoBack bookmark right after the final sentence.
# -----------------------------------------------------------------------

# Helper: appends a brand-new, perfectly empty paragraph at the end of the
# document. A placeholder character is inserted together with the new
# paragraph mark and removed straight afterwards - inserting the mark on
# its own tends to leave a stray empty <w:r/> behind.
function New-TrailingEmptyParagraph($doc) {
    $count = $doc.Paragraphs.Count
    $tail = $doc.Paragraphs.Item($count).Range
    $tail.Collapse(0)
    $tail.InsertAfter([char]13 + "X")
    $newTail = $doc.Paragraphs.Item($doc.Paragraphs.Count).Range
    $doc.Range($newTail.Start, $newTail.Start + 1).Delete()
}

# Create all four trailing paragraphs while the "current" paragraph
# formatting is still plain, then fill in the text afterwards - this
# avoids the bold/underline of the date line leaking into the paragraphs
# that come after it.
New-TrailingEmptyParagraph $d
New-TrailingEmptyParagraph $d
New-TrailingEmptyParagraph $d
New-TrailingEmptyParagraph $d

$total = $d.Paragraphs.Count

# Paragraph (total-2) becomes the bold/underlined date line.
$dateParaIndex = $total - 2
$dateRange = $d.Paragraphs.Item($dateParaIndex).Range
$dateRange.InsertBefore("12/06/2018")
$dateRange = $d.Paragraphs.Item($dateParaIndex).Range
$dateRange.Font.Bold = 1
$dateRange.Font.Underline = 1

# Last paragraph becomes the closing sentence. An extra placeholder is
# appended after the real text so that the position right after the
# sentence is never the literal end of the document/run while the
# bookmark is created there; the placeholder is stripped away right
# after (doing this avoids the bookmark ending up mis-anchored).
$lastParaIndex = $total
$lastRange = $d.Paragraphs.Item($lastParaIndex).Range
$lastRange.InsertBefore("Avancement dans le rapport.ZZZ")

$searchRange = $d.Content
$searchFind = $searchRange.Find
$searchFind.ClearFormatting()
$searchFind.Execute(
    "Avancement dans le rapport.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$searchRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $searchRange)

$placeholderStart = $searchRange.Start
$d.Range($placeholderStart, $placeholderStart + 3).Delete()
